$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Fuels")

# Add moisture content data (0) for sinter plant (coke) and blast furnace (PCI coal)
$ws.Range("E10").Value = 0
$ws.Range("E12").Value = 0

# Update the selected/active cell on the sheet
$ws.Range("D12").Select()
